$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.05088209513382935
$ws.Range("D2").Value = 0.09244243276017983
$ws.Range("E2").Value = 0.4343060183838787
$ws.Range("F2").Value = 2.285403960572268
$ws.Range("G2").Value = 2.564910388463545
$ws.Range("H2").Value = 1.534877293481372
$ws.Range("C3").Value = 0.04499111415226764
$ws.Range("D3").Value = 0.08099312179224682
$ws.Range("E3").Value = 0.3774527995275889
$ws.Range("F3").Value = 2.053156900224678
$ws.Range("G3").Value = 2.275893840270044
$ws.Range("H3").Value = 1.407448340504004
$ws.Range("C4").Value = 0.04140688078128107
$ws.Range("D4").Value = 0.07402400667483278
$ws.Range("E4").Value = 0.3428124898187264
$ws.Range("F4").Value = 1.912316711634247
$ws.Range("G4").Value = 2.100249222090611
$ws.Range("H4").Value = 1.330500229916026
$ws.Range("C5").Value = 0.03995402274648541
$ws.Range("D5").Value = 0.07119808863352262
$ws.Range("E5").Value = 0.3287567351206917
$ws.Range("F5").Value = 1.855343698008113
$ws.Range("G5").Value = 2.02910102924011
$ws.Range("H5").Value = 1.299455806609899
$ws.Range("C6").Value = 0.03971322982907566
$ws.Range("D6").Value = 0.07072966030467853
$ws.Range("E6").Value = 0.3264262411716601
$ws.Range("F6").Value = 1.845908133199885
$ws.Range("G6").Value = 2.017312013859168
$ws.Range("H6").Value = 1.294319414592678
$ws.Range("C7").Value = 0.04138725633642082
$ws.Range("D7").Value = 0.07398584006757858
$ws.Range("E7").Value = 0.3426226931926806
$ws.Range("F7").Value = 1.911546679905854
$ws.Range("G7").Value = 2.099287993742507
$ws.Range("H7").Value = 1.330080305065678
$ws.Range("C8").Value = 0.048843750634461
$ws.Range("D8").Value = 0.08848133144697101
$ws.Range("E8").Value = 0.4146432302580934
$ws.Range("F8").Value = 2.204946202125285
$ws.Range("G8").Value = 2.46486445801628
$ws.Range("H8").Value = 1.490663432560325
$ws.Range("C9").Value = 0.06375108278027142
$ws.Range("D9").Value = 0.117447993817521
$ws.Range("E9").Value = 0.5583251547865729
$ws.Range("F9").Value = 2.79532394359515
$ws.Range("G9").Value = 3.197443359351439
$ws.Range("H9").Value = 1.816439590657353
$ws.Range("C10").Value = 0.07491319404149976
$ws.Range("D10").Value = 0.1391449197427335
$ws.Range("E10").Value = 0.6658540739310013
$ws.Range("F10").Value = 3.239773343520341
$ws.Range("G10").Value = 3.747132436533263
$ws.Range("H10").Value = 2.063304805149755
$ws.Range("C11").Value = 0.08004481728909241
$ws.Range("D11").Value = 0.1491249759868936
$ws.Range("E11").Value = 0.7153062832119872
$ws.Range("F11").Value = 3.444639034242755
$ws.Range("G11").Value = 4.00012231043263
$ws.Range("H11").Value = 2.177445210765768
$ws.Range("C12").Value = 0.0819964595482503
$ws.Range("D12").Value = 0.1529216098639381
$ws.Range("E12").Value = 0.7341187954332895
$ws.Range("F12").Value = 3.522631409936196
$ws.Range("G12").Value = 4.096380848006049
$ws.Range("H12").Value = 2.22094871134567
$ws.Range("C13").Value = 0.08157575385052951
$ws.Range("D13").Value = 0.1521031369487673
$ws.Range("E13").Value = 0.73006321326136
$ws.Range("F13").Value = 3.505815472973211
$ws.Range("G13").Value = 4.075628964692498
$ws.Range("H13").Value = 2.211566688482378
$ws.Range("C14").Value = 0.08020520773034434
$ws.Range("D14").Value = 0.1494369692640305
$ws.Range("E14").Value = 0.7168522236950423
$ws.Range("F14").Value = 3.451047047003783
$ws.Range("G14").Value = 4.008032199225568
$ws.Range("H14").Value = 2.181018542013135
$ws.Range("C15").Value = 0.07936682344221424
$ws.Range("D15").Value = 0.147806179020904
$ws.Range("E15").Value = 0.708771583954217
$ws.Range("F15").Value = 3.417554602524831
$ws.Range("G15").Value = 3.96668777174574
$ws.Range("H15").Value = 2.162344012835035
$ws.Range("C16").Value = 0.07457897080197995
$ws.Range("D16").Value = 0.1384950473253639
$ws.Range("E16").Value = 0.6626338051535754
$ws.Range("F16").Value = 3.226441281925361
$ws.Range("G16").Value = 3.73066093909091
$ws.Range("H16").Value = 2.055883877073825
$ws.Range("C17").Value = 0.07165604684527693
$ws.Range("D17").Value = 0.1328122544848327
$ws.Range("E17").Value = 0.6344734356464699
$ws.Range("F17").Value = 3.109906547216184
$ws.Range("G17").Value = 3.586641909769185
$ws.Range("H17").Value = 1.991056810500083
$ws.Range("C18").Value = 0.0699798932133433
$ws.Range("D18").Value = 0.1295539192347519
$ws.Range("E18").Value = 0.6183263086667239
$ws.Range("F18").Value = 3.043129669826271
$ws.Range("G18").Value = 3.504079926245197
$ws.Range("H18").Value = 1.95394218506118
$ws.Range("C19").Value = 0.06941322161699759
$ws.Range("D19").Value = 0.1284524172932606
$ws.Range("E19").Value = 0.6128674879011555
$ws.Range("F19").Value = 3.020562399879083
$ws.Range("G19").Value = 3.476171866273489
$ws.Range("H19").Value = 1.941404895010123
$ws.Range("C20").Value = 0.0719666713361562
$ws.Range("D20").Value = 0.1334161246627446
$ws.Range("E20").Value = 0.63746592009241
$ws.Range("F20").Value = 3.122285708113651
$ws.Range("G20").Value = 3.601944379440454
$ws.Range("H20").Value = 1.997939824818161
$ws.Range("C21").Value = 0.08060753682967459
$ws.Range("D21").Value = 0.1502196021049826
$ws.Range("E21").Value = 0.7207302043796489
$ws.Range("F21").Value = 3.467122389459746
$ws.Range("G21").Value = 4.027874326515303
$ws.Range("H21").Value = 2.189983512906053
$ws.Range("C22").Value = 0.08630417082480335
$ws.Range("D22").Value = 0.1613038644017593
$ws.Range("E22").Value = 0.7756540621339951
$ws.Range("F22").Value = 3.694921204907018
$ws.Range("G22").Value = 4.308923702799973
$ws.Range("H22").Value = 2.317141151680062
$ws.Range("C23").Value = 0.08325903397782497
$ws.Range("D23").Value = 0.1553780927343666
$ws.Range("E23").Value = 0.7462908648888202
$ws.Range("F23").Value = 3.573108878531798
$ws.Range("G23").Value = 4.158665295903404
$ws.Range("H23").Value = 2.249118535815228
$ws.Range("C24").Value = 0.07182622473997924
$ws.Range("D24").Value = 0.1331430877323783
$ws.Range("E24").Value = 0.6361128857112135
$ws.Range("F24").Value = 3.116688404811185
$ws.Range("G24").Value = 3.595025398961468
$ws.Range("H24").Value = 1.994827531493968
$ws.Range("C25").Value = 0.05968409346859005
$ws.Range("D25").Value = 0.1095451899372506
$ws.Range("E25").Value = 0.5191453392989445
$ws.Range("F25").Value = 2.633846868522511
$ws.Range("G25").Value = 2.997389641847974
$ws.Range("H25").Value = 1.727054494182937
